# This script reproduces the "Natmi following Dr Hou advice" edit:
# a new sending-cluster group ("ECs") is added to the Has2-Cd44 ligand/receptor
# table, so the sending x target cluster grid grows from 2x3 to 3x3 (rows 2-10),
# and every NATMI statistic column (E:T) is recomputed for the new 3-cluster grid.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (Has2 -> Cd44)
$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Has2"
$ws.Range("C2").Value2 = "Cd44"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 0.6666666666666666
$ws.Range("G2").Value2 = 0.3897156666666666
$ws.Range("H2").Value2 = 1.169147
$ws.Range("I2").Value2 = 0.01604749516233658
$ws.Range("J2").Value2 = 0.01604749516233658
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 297.8183156666666
$ws.Range("N2").Value2 = 893.454947
$ws.Range("O2").Value2 = 0.8852156413092672
$ws.Range("P2").Value2 = 0.8852156413092673
$ws.Range("Q2").Value2 = 116.0644634355788
$ws.Range("R2").Value2 = 1044.580170920209
$ws.Range("S2").Value2 = 0.01420549372153514
$ws.Range("T2").Value2 = 0.01420549372153514

# Row 3: ECs -> FAPs (Has2 -> Cd44)
$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "Has2"
$ws.Range("C3").Value2 = "Cd44"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 2
$ws.Range("F3").Value2 = 0.6666666666666666
$ws.Range("G3").Value2 = 0.3897156666666666
$ws.Range("H3").Value2 = 1.169147
$ws.Range("I3").Value2 = 0.01604749516233658
$ws.Range("J3").Value2 = 0.01604749516233658
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 24.34034433333333
$ws.Range("N3").Value2 = 73.021033
$ws.Range("O3").Value2 = 0.07234764413494278
$ws.Range("P3").Value2 = 0.0723476441349428
$ws.Range("Q3").Value2 = 9.485813518761221
$ws.Range("R3").Value2 = 85.372321668851
$ws.Range("S3").Value2 = 0.001160998469261943
$ws.Range("T3").Value2 = 0.001160998469261943

# Row 4: ECs -> sCs (Has2 -> Cd44)
$ws.Range("A4").Value2 = "ECs"
$ws.Range("B4").Value2 = "Has2"
$ws.Range("C4").Value2 = "Cd44"
$ws.Range("D4").Value2 = "sCs"
$ws.Range("E4").Value2 = 2
$ws.Range("F4").Value2 = 0.6666666666666666
$ws.Range("G4").Value2 = 0.3897156666666666
$ws.Range("H4").Value2 = 1.169147
$ws.Range("I4").Value2 = 0.01604749516233658
$ws.Range("J4").Value2 = 0.01604749516233658
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 14.277234
$ws.Range("N4").Value2 = 42.831702
$ws.Range("O4").Value2 = 0.04243671455578994
$ws.Range("P4").Value2 = 0.04243671455578994
$ws.Range("Q4").Value2 = 5.564061766466
$ws.Range("R4").Value2 = 50.07655589819399
$ws.Range("S4").Value2 = 0.0006810029715394973
$ws.Range("T4").Value2 = 0.0006810029715394976

# Row 5: FAPs -> ECs (Has2 -> Cd44)
$ws.Range("A5").Value2 = "FAPs"
$ws.Range("B5").Value2 = "Has2"
$ws.Range("C5").Value2 = "Cd44"
$ws.Range("D5").Value2 = "ECs"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 23.535792
$ws.Range("H5").Value2 = 70.607376
$ws.Range("I5").Value2 = 0.9691437644584301
$ws.Range("J5").Value2 = 0.9691437644584302
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 297.8183156666666
$ws.Range("N5").Value2 = 893.454947
$ws.Range("O5").Value2 = 0.8852156413092672
$ws.Range("P5").Value2 = 0.8852156413092673
$ws.Range("Q5").Value2 = 7009.389931321008
$ws.Range("R5").Value2 = 63084.50938188907
$ws.Range("S5").Value2 = 0.8579012189759466
$ws.Range("T5").Value2 = 0.8579012189759468

# Row 6: FAPs -> FAPs (Has2 -> Cd44)
$ws.Range("A6").Value2 = "FAPs"
$ws.Range("B6").Value2 = "Has2"
$ws.Range("C6").Value2 = "Cd44"
$ws.Range("D6").Value2 = "FAPs"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 23.535792
$ws.Range("H6").Value2 = 70.607376
$ws.Range("I6").Value2 = 0.9691437644584301
$ws.Range("J6").Value2 = 0.9691437644584302
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 24.34034433333333
$ws.Range("N6").Value2 = 73.021033
$ws.Range("O6").Value2 = 0.07234764413494278
$ws.Range("P6").Value2 = 0.0723476441349428
$ws.Range("Q6").Value2 = 572.8692814377121
$ws.Range("R6").Value2 = 5155.823532939408
$ws.Range("S6").Value2 = 0.07011526818663731
$ws.Range("T6").Value2 = 0.07011526818663733

# Row 7: FAPs -> sCs (Has2 -> Cd44)
$ws.Range("A7").Value2 = "FAPs"
$ws.Range("B7").Value2 = "Has2"
$ws.Range("C7").Value2 = "Cd44"
$ws.Range("D7").Value2 = "sCs"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 23.535792
$ws.Range("H7").Value2 = 70.607376
$ws.Range("I7").Value2 = 0.9691437644584301
$ws.Range("J7").Value2 = 0.9691437644584302
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 14.277234
$ws.Range("N7").Value2 = 42.831702
$ws.Range("O7").Value2 = 0.04243671455578994
$ws.Range("P7").Value2 = 0.04243671455578994
$ws.Range("Q7").Value2 = 336.026009759328
$ws.Range("R7").Value2 = 3024.234087833952
$ws.Range("S7").Value2 = 0.04112727729584612
$ws.Range("T7").Value2 = 0.04112727729584612

# Row 8: sCs -> ECs (Has2 -> Cd44)
$ws.Range("A8").Value2 = "sCs"
$ws.Range("B8").Value2 = "Has2"
$ws.Range("C8").Value2 = "Cd44"
$ws.Range("D8").Value2 = "ECs"
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 0.3596323333333333
$ws.Range("H8").Value2 = 1.078897
$ws.Range("I8").Value2 = 0.01480874037923328
$ws.Range("J8").Value2 = 0.01480874037923328
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 297.8183156666666
$ws.Range("N8").Value2 = 893.454947
$ws.Range("O8").Value2 = 0.8852156413092672
$ws.Range("P8").Value2 = 0.8852156413092673
$ws.Range("Q8").Value2 = 107.1050957726065
$ws.Range("R8").Value2 = 963.9458619534589
$ws.Range("S8").Value2 = 0.01310892861178543
$ws.Range("T8").Value2 = 0.01310892861178543

# Row 9: sCs -> FAPs (Has2 -> Cd44)
$ws.Range("A9").Value2 = "sCs"
$ws.Range("B9").Value2 = "Has2"
$ws.Range("C9").Value2 = "Cd44"
$ws.Range("D9").Value2 = "FAPs"
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 0.3596323333333333
$ws.Range("H9").Value2 = 1.078897
$ws.Range("I9").Value2 = 0.01480874037923328
$ws.Range("J9").Value2 = 0.01480874037923328
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 24.34034433333333
$ws.Range("N9").Value2 = 73.021033
$ws.Range("O9").Value2 = 0.07234764413494278
$ws.Range("P9").Value2 = 0.0723476441349428
$ws.Range("Q9").Value2 = 8.753574826733445
$ws.Range("R9").Value2 = 78.782173440601
$ws.Range("S9").Value2 = 0.001071377479043527
$ws.Range("T9").Value2 = 0.001071377479043527

# Row 10: sCs -> sCs (Has2 -> Cd44)
$ws.Range("A10").Value2 = "sCs"
$ws.Range("B10").Value2 = "Has2"
$ws.Range("C10").Value2 = "Cd44"
$ws.Range("D10").Value2 = "sCs"
$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 0.3596323333333333
$ws.Range("H10").Value2 = 1.078897
$ws.Range("I10").Value2 = 0.01480874037923328
$ws.Range("J10").Value2 = 0.01480874037923328
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 14.277234
$ws.Range("N10").Value2 = 42.831702
$ws.Range("O10").Value2 = 0.04243671455578994
$ws.Range("P10").Value2 = 0.04243671455578994
$ws.Range("Q10").Value2 = 5.134554976966
$ws.Range("R10").Value2 = 46.210994792694
$ws.Range("S10").Value2 = 0.0006284342884043232
$ws.Range("T10").Value2 = 0.0006284342884043233
